# Add two new worksheets ("SendEmail" and "ReportPDF") right after the
# "Login" sheet, matching the commit "Added email and pdf report".

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")

# Insert SendEmail right after Login, then ReportPDF right after SendEmail.
$sendEmail = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$sendEmail.Name = "SendEmail"

$reportPdf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sendEmail)
$reportPdf.Name = "ReportPDF"

# ---------------------------------------------------------------------
# SendEmail sheet
# ---------------------------------------------------------------------

# Header row (columns A-L, N first; M is filled in later, after the data
# row, to reproduce the shared-string insertion order of the authored
# workbook).
$sendEmail.Range("A1").Value = "Subject"
$sendEmail.Range("B1").Value = "Body"
$sendEmail.Range("C1").Value = "To"
$sendEmail.Range("D1").Value = "From"
$sendEmail.Range("E1").Value = "ServerHostName"
$sendEmail.Range("F1").Value = "ServerPort"
$sendEmail.Range("G1").Value = "UseSSL"
$sendEmail.Range("H1").Value = "Username"
$sendEmail.Range("I1").Value = "Password"
$sendEmail.Range("J1").Value = "SendEmailOnFailure"
$sendEmail.Range("K1").Value = "SendEmailOnSuccess"
$sendEmail.Range("L1").Value = "SendZippedReportOnComplete"
$sendEmail.Range("N1").Value = "PDFReportCustomStyleSheet"

# Data row
$sendEmail.Range("A2").Value = "Build AA Automation 2019"

# Body (B2) is left blank but wraps text, like the source workbook.
$sendEmail.Range("B2").WrapText = $true

$sendEmail.Range("C2").Value = "mmargasagayam@abacusnext.com"
$sendEmail.Range("D2").Value = "amicustestmk1@gmail.com"
$sendEmail.Range("E2").Value = "smtp.gmail.com"
$sendEmail.Range("F2").Value = 587
$sendEmail.Range("G2").Value = $true
$sendEmail.Range("H2").Value = "amicustestmk1@gmail.com"
$sendEmail.Range("I2").Value = "0nXTeam123$$"
$sendEmail.Range("J2").Value = $false
$sendEmail.Range("K2").Value = $false
$sendEmail.Range("L2").Value = $false
$sendEmail.Range("M2").Value = $true

# M1 header is typed in last (matches authored shared-string order).
$sendEmail.Range("M1").Value = "SendPDFReportOnComplete"

# Hyperlinks (Excel auto-links e-mail addresses; recreate the same effect).
$sendEmail.Hyperlinks.Add($sendEmail.Range("C2"), "mailto:mmargasagayam@abacusnext.com") | Out-Null
$sendEmail.Hyperlinks.Add($sendEmail.Range("D2"), "mailto:amicustestmk1@gmail.com") | Out-Null
$sendEmail.Hyperlinks.Add($sendEmail.Range("H2"), "mailto:amicustestmk1@gmail.com") | Out-Null

# Approximate column widths (best-fit sizing from the authored workbook).
$sendEmail.Columns.Item(1).ColumnWidth = 23.43
$sendEmail.Columns.Item(2).ColumnWidth = 69.02
$sendEmail.Columns.Item(3).ColumnWidth = 31.74
$sendEmail.Columns.Item(4).ColumnWidth = 25.02
$sendEmail.Columns.Item(5).ColumnWidth = 15.31
$sendEmail.Columns.Item(6).ColumnWidth = 9.59
$sendEmail.Columns.Item(7).ColumnWidth = 6.31
$sendEmail.Columns.Item(8).ColumnWidth = 25.02
$sendEmail.Columns.Item(9).ColumnWidth = 13.45
$sendEmail.Columns.Item(10).ColumnWidth = 18.45
$sendEmail.Columns.Item(11).ColumnWidth = 19.02
$sendEmail.Columns.Item(12).ColumnWidth = 28.59
$sendEmail.Columns.Item(13).ColumnWidth = 25.88
$sendEmail.Columns.Item(14).ColumnWidth = 26.45

# ---------------------------------------------------------------------
# ReportPDF sheet
# ---------------------------------------------------------------------

$reportPdf.Range("A1").Value = "PdfName"
$reportPdf.Range("B1").Value = "PdfDirectoryPath"
$reportPdf.Range("C1").Value = "Xml"
$reportPdf.Range("D1").Value = "Details"
$reportPdf.Range("E1").Value = "DeleteRanorexReport"

$reportPdf.Columns.Item(1).ColumnWidth = 8.45
$reportPdf.Columns.Item(2).ColumnWidth = 15.45
$reportPdf.Columns.Item(3).ColumnWidth = 3.59
$reportPdf.Columns.Item(4).ColumnWidth = 6.31
$reportPdf.Columns.Item(5).ColumnWidth = 19.74

$reportPdf.Range("E1").Select() | Out-Null

# ---------------------------------------------------------------------
# Selection bookkeeping: user ends up back on Login (cell H5) before
# finally leaving SendEmail (cell B2) as the active sheet/cell.
# ---------------------------------------------------------------------

$loginSheet.Select()
$loginSheet.Range("H5").Select() | Out-Null

$sendEmail.Select()
$sendEmail.Range("B2").Select() | Out-Null
